# Update the "through 11-13" -> "through 11-14" dataset (one more day of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet (tab) name reflects the date the data run through.
$ws.Name = "Through 2022-11-14"

# Update the "November (through 11-13)" label used in column A / shared strings.
$ws.Cells.Replace("through 11-13", "through 11-14")

# Row 12 = November monthly counts for years 2015-2022 (columns B-I).
$ws.Cells.Item(12, 2).Value = 17
$ws.Cells.Item(12, 3).Value = 35
$ws.Cells.Item(12, 4).Value = 62
$ws.Cells.Item(12, 5).Value = 29
$ws.Cells.Item(12, 6).Value = 23
$ws.Cells.Item(12, 7).Value = 87
$ws.Cells.Item(12, 8).Value = 98
$ws.Cells.Item(12, 9).Value = 47

# Row 13 = Total counts for years 2015-2022 (columns B-I), updated to include
# the newly added November data.
$ws.Cells.Item(13, 2).Value = 275
$ws.Cells.Item(13, 3).Value = 521
$ws.Cells.Item(13, 4).Value = 772
$ws.Cells.Item(13, 5).Value = 644
$ws.Cells.Item(13, 6).Value = 505
$ws.Cells.Item(13, 7).Value = 1144
$ws.Cells.Item(13, 8).Value = 1539
$ws.Cells.Item(13, 9).Value = 1445
